$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 already has B/C formatted as time (style index 2 -> numFmt 20) and
# D as wrap-text string (style index 7). Copy that formatting down to row 34
# so the new row matches the existing table's look, then fill in the values.
$ws.Range("B33:D33").Copy()
$ws.Range("B34:D34").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B34").Value = 0.10416666666666667
$ws.Range("C34").Value = 0.16666666666666666
$ws.Range("D34").Value = "There are some bugs with the IK interpolation while movin."

$ws.Range("C35").Select()
